# Insert two new data rows at the top of the data block (rows 36-37),
# pushing the existing rows (old 36..126) down to (38..128).
# Then populate the two new rows with their data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("36:37").Insert()

# Row 36 (new)
$ws.Range("A36").Value = 4
$ws.Range("B36").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C36").Value = 'Los Lagos'
$ws.Range("D36").Value = 44414
$ws.Range("E36").Value = 10
$ws.Range("F36").Value = 100112021
$ws.Range("G36").Value = 'Ají'
$ws.Range("H36").Value = 'Inferno'
$ws.Range("I36").Value = 'Primera'
$ws.Range("J36").Value = 80
$ws.Range("K36").Value = 35000
$ws.Range("L36").Value = 35000
$ws.Range("M36").Value = 35000
$ws.Range("N36").Value = '$/caja 12 kilos'
$ws.Range("O36").Value = 'Región de Arica y Parinacota'
$ws.Range("P36").Value = 2917
$ws.Range("Q36").Value = 12
$ws.Range("R36").Value = 'Hortaliza'

# Row 37 (new)
$ws.Range("A37").Value = 4
$ws.Range("B37").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C37").Value = 'Los Lagos'
$ws.Range("D37").Value = 44414
$ws.Range("E37").Value = 10
$ws.Range("F37").Value = 100112021
$ws.Range("G37").Value = 'Ají'
$ws.Range("H37").Value = 'Inferno'
$ws.Range("I37").Value = 'Segunda'
$ws.Range("J37").Value = 80
$ws.Range("K37").Value = 28000
$ws.Range("L37").Value = 28000
$ws.Range("M37").Value = 28000
$ws.Range("N37").Value = '$/caja 12 kilos'
$ws.Range("O37").Value = 'Región de Arica y Parinacota'
$ws.Range("P37").Value = 2333
$ws.Range("Q37").Value = 12
$ws.Range("R37").Value = 'Hortaliza'
